# Auto-generated edit script applying diff to cryptos.xlsx
# Updates Price (D) and Volume(1h) (E) columns for several rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.406.71'
$ws.Range('D2').Style = $ws.Range('B2').Style
$ws.Range('E2').Value = '  -0.47%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.609.97'
$ws.Range('D3').Style = $ws.Range('B3').Style
$ws.Range('E3').Value = '  +0.28%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '583.12'
$ws.Range('D5').Style = $ws.Range('B5').Style
$ws.Range('E5').Value = '  +2.43%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.12'
$ws.Range('D6').Style = $ws.Range('B6').Style
$ws.Range('E6').Value = '  +0.25%  '
$ws.Range('E8').Value = '  -0.69%  '
$ws.Range('E9').Value = '  +0.38%  '
$ws.Range('E10').Value = '  -1.37%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.155'
$ws.Range('D11').Style = $ws.Range('B11').Style
$ws.Range('E11').Value = '  +0.61%  '
$ws.Range('E12').Value = '  +1.22%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.071.32'
$ws.Range('D13').Style = $ws.Range('B13').Style
$ws.Range('E13').Value = '  +0.10%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '24.69'
$ws.Range('D14').Style = $ws.Range('B14').Style
$ws.Range('E14').Value = '  +4.67%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '60.406.01'
$ws.Range('D15').Style = $ws.Range('B15').Style
$ws.Range('E15').Value = '  -0.51%  '
$ws.Range('E16').Value = '  -0.08%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.613.71'
$ws.Range('D17').Style = $ws.Range('B17').Style
$ws.Range('E17').Value = '  +0.00%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.35'
$ws.Range('D18').Style = $ws.Range('B18').Style
$ws.Range('E18').Value = '  +0.51%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.63'
$ws.Range('D19').Style = $ws.Range('B19').Style
$ws.Range('E19').Value = '  -1.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '346.38'
$ws.Range('D20').Style = $ws.Range('B20').Style
$ws.Range('E20').Value = '  -0.10%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.91'
$ws.Range('D21').Style = $ws.Range('B21').Style
$ws.Range('E21').Value = '  -2.91%  '
$ws.Range('E22').Value = '  -0.16%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.535'
$ws.Range('D23').Style = $ws.Range('B23').Style
$ws.Range('E23').Value = '  +2.99%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.67'
$ws.Range('D24').Style = $ws.Range('B24').Style
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.999'
$ws.Range('D25').Style = $ws.Range('B25').Style
$ws.Range('E25').Value = '  +0.44%  '
$ws.Range('E26').Value = '  +0.20%  '
$ws.Range('E27').Value = '  +3.14%  '
$ws.Range('E28').Value = '  +4.87%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0798'
$ws.Range('D29').Style = $ws.Range('B29').Style
$ws.Range('E29').Value = '  +0.46%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '169.11'
$ws.Range('D30').Style = $ws.Range('B30').Style
$ws.Range('E30').Value = '  +4.70%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.43'
$ws.Range('D31').Style = $ws.Range('B31').Style
$ws.Range('E31').Value = '  +2.29%  '
$ws.Range('E32').Value = '  +0.11%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.46'
$ws.Range('D33').Style = $ws.Range('B33').Style
$ws.Range('E33').Value = '  -0.11%  '
$ws.Range('E34').Value = '  +9.12%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.28'
$ws.Range('D35').Style = $ws.Range('B35').Style
$ws.Range('E35').Value = '  +0.35%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.01'
$ws.Range('D36').Style = $ws.Range('B36').Style
$ws.Range('E36').Value = '  +5.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.65'
$ws.Range('D37').Style = $ws.Range('B37').Style
$ws.Range('E37').Value = '  +4.11%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '318.83'
$ws.Range('D38').Style = $ws.Range('B38').Style
$ws.Range('E38').Value = '  +7.52%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '38.34'
$ws.Range('D39').Style = $ws.Range('B39').Style
$ws.Range('E39').Value = '  +1.64%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.90'
$ws.Range('D40').Style = $ws.Range('B40').Style
$ws.Range('E40').Value = '  +2.68%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.850'
$ws.Range('D41').Style = $ws.Range('B41').Style
$ws.Range('E41').Value = '  -0.79%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '135.56'
$ws.Range('D42').Style = $ws.Range('B42').Style
$ws.Range('E42').Value = '  -2.26%  '
$ws.Range('E43').Value = '  +1.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.998'
$ws.Range('D44').Style = $ws.Range('B44').Style
$ws.Range('E44').Value = '  +0.30%  '
$ws.Range('E45').Value = '  +1.64%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.609'
$ws.Range('D46').Style = $ws.Range('B46').Style
$ws.Range('E46').Value = '  +0.75%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.02'
$ws.Range('D47').Style = $ws.Range('B47').Style
$ws.Range('E47').Value = '  +4.48%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0550'
$ws.Range('D48').Style = $ws.Range('B48').Style
$ws.Range('E48').Value = '  -0.42%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '20.06'
$ws.Range('D49').Style = $ws.Range('B49').Style
$ws.Range('E49').Value = '  +1.71%  '
$ws.Range('E50').Value = '  +0.18%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '10.76'
$ws.Range('D51').Style = $ws.Range('B51').Style
$ws.Range('E51').Value = '  +0.60%  '
